$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 39657
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 39657
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 39657
$ws.Range("N3").Value = -39885
$ws.Range("H11").Value = 87.333336
$ws.Range("I11").Value = 87.333336
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 87.333336
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 52.666664
$ws.Range("H54").Value = 10250
$ws.Range("I54").Value = 10250
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 10250
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -9764
$ws.Range("H74").Value = 3312.5
$ws.Range("I74").Value = 3583.3333
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 3583.3333
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -2647.3333
$ws.Range("N74").Value = -4372
$ws.Range("H77").Value = 3312.5
$ws.Range("I77").Value = 3583.3333
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 17916.6665
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -13236.6665
$ws.Range("N77").Value = -21860
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H102").Value = 39657
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 39657
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 39657
$ws.Range("N102").Value = -46147
$ws.Range("H104").Value = 222.14285
$ws.Range("I104").Value = 222.14285
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 666.4285500000001
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 1080.57145
$ws.Range("H110").Value = 33314
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 33314
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 33314
$ws.Range("N110").Value = -41494
$ws.Range("H132").Value = 3751.7568
$ws.Range("I132").Value = 3359.0605
$ws.Range("J132").Value = 6991.5
$ws.Range("K132").Value = 10077.1815
$ws.Range("L132").Value = 20974.5
$ws.Range("M132").Value = -7547.181500000001
$ws.Range("N132").Value = -26034.5
$ws.Range("H137").Value = 2755.4443
$ws.Range("I137").Value = 2141.5
$ws.Range("J137").Value = 3983.3333
$ws.Range("K137").Value = 6424.5
$ws.Range("L137").Value = 11949.9999
$ws.Range("M137").Value = -3874.5
$ws.Range("N137").Value = -17049.9999
$ws.Range("H141").Value = 4200
$ws.Range("I141").Value = 4200
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 12600
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -7420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2093.0557
$ws.Range("I2").Value = 1346.8334
$ws.Range("J2").Value = 3585.5
$ws.Range("K2").Value = 1346.8334
$ws.Range("L2").Value = 3585.5
$ws.Range("M2").Value = -1233.8334
$ws.Range("N2").Value = -3811.5
$ws.Range("H5").Value = 42.11111
$ws.Range("I5").Value = 42.285713
$ws.Range("J5").Value = 41.5
$ws.Range("K5").Value = 42.285713
$ws.Range("L5").Value = 41.5
$ws.Range("M5").Value = 69.714287
$ws.Range("N5").Value = -265.5
$ws.Range("H32").Value = 25360.967
$ws.Range("I32").Value = 17822.074
$ws.Range("J32").Value = 76248.5
$ws.Range("K32").Value = 17822.074
$ws.Range("L32").Value = 76248.5
$ws.Range("M32").Value = -17535.074
$ws.Range("N32").Value = -76822.5
$ws.Range("H37").Value = 23330.555
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 23330.555
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 23330.555
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -23876.555
$ws.Range("H61").Value = 3984.625
$ws.Range("I61").Value = 2997.8
$ws.Range("J61").Value = 5629.3335
$ws.Range("K61").Value = 2997.8
$ws.Range("L61").Value = 5629.3335
$ws.Range("M61").Value = -2785.8
$ws.Range("N61").Value = -6053.3335
$ws.Range("H74").Value = 2199.2
$ws.Range("I74").Value = 2110.7778
$ws.Range("J74").Value = 2995
$ws.Range("K74").Value = 2110.7778
$ws.Range("L74").Value = 2995
$ws.Range("M74").Value = -1236.7778
$ws.Range("N74").Value = -4743
$ws.Range("H77").Value = 2199.2
$ws.Range("I77").Value = 2110.7778
$ws.Range("J77").Value = 2995
$ws.Range("K77").Value = 10553.889
$ws.Range("L77").Value = 14975
$ws.Range("M77").Value = -6185.888999999999
$ws.Range("N77").Value = -23711
$ws.Range("H80").Value = 39998.57
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 39998.57
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 39998.57
$ws.Range("N80").Value = -41994.57
$ws.Range("H83").Value = 39998.57
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 39998.57
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 119995.71
$ws.Range("N83").Value = -129979.71
$ws.Range("H88").Value = 9948.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9948.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 9948.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -10760.5
$ws.Range("H91").Value = 9948.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9948.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 9948.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -12756.5
$ws.Range("H116").Value = 2093.0557
$ws.Range("I116").Value = 1346.8334
$ws.Range("J116").Value = 3585.5
$ws.Range("K116").Value = 1346.8334
$ws.Range("L116").Value = 3585.5
$ws.Range("M116").Value = 947.1666
$ws.Range("N116").Value = -8173.5
$ws.Range("H136").Value = 3984.625
$ws.Range("I136").Value = 2997.8
$ws.Range("J136").Value = 5629.3335
$ws.Range("K136").Value = 8993.400000000001
$ws.Range("L136").Value = 16888.0005
$ws.Range("M136").Value = -6443.400000000001
$ws.Range("N136").Value = -21988.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2093.0557
$ws.Range("I3").Value = 1346.8334
$ws.Range("J3").Value = 3585.5
$ws.Range("K3").Value = 1346.8334
$ws.Range("L3").Value = 3585.5
$ws.Range("M3").Value = -1232.8334
$ws.Range("N3").Value = -3813.5
$ws.Range("H4").Value = 42.11111
$ws.Range("I4").Value = 42.285713
$ws.Range("J4").Value = 41.5
$ws.Range("K4").Value = 42.285713
$ws.Range("L4").Value = 41.5
$ws.Range("M4").Value = 72.714287
$ws.Range("N4").Value = -271.5
$ws.Range("H9").Value = 10000
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 10000
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10336
$ws.Range("H94").Value = 1474.24
$ws.Range("I94").Value = 1327.3334
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 1327.3334
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -876.3334
$ws.Range("N94").Value = -5902
$ws.Range("H95").Value = 18574.666
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 18574.666
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 18574.666
$ws.Range("N95").Value = -24066.666
$ws.Range("H100").Value = 49499.5
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 49499.5
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 49499.5
$ws.Range("N100").Value = -51663.5
$ws.Range("H103").Value = 21830
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 21830
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 21830
$ws.Range("N103").Value = -24174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6385.1665
$ws.Range("I31").Value = 5332
$ws.Range("J31").Value = 7438.3335
$ws.Range("K31").Value = 5332
$ws.Range("L31").Value = 7438.3335
$ws.Range("M31").Value = -5037
$ws.Range("N31").Value = -8028.3335
$ws.Range("H34").Value = 6385.1665
$ws.Range("I34").Value = 5332
$ws.Range("J34").Value = 7438.3335
$ws.Range("K34").Value = 5332
$ws.Range("L34").Value = 7438.3335
$ws.Range("M34").Value = -5130
$ws.Range("N34").Value = -7842.3335
$ws.Range("H43").Value = 50000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 50000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50368
$ws.Range("H62").Value = 10001.25
$ws.Range("I62").Value = 10001.667
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 10001.667
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -9377.666999999999
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 10001.25
$ws.Range("I65").Value = 10001.667
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 50008.335
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -46888.335
$ws.Range("N65").Value = -56240
$ws.Range("H86").Value = 3533.3333
$ws.Range("I86").Value = 3533.3333
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3533.3333
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2410.3333
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 3533.3333
$ws.Range("I89").Value = 3533.3333
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 17666.6665
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -12050.6665
$ws.Range("N89").ClearContents()
$ws.Range("H101").Value = 50000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 50000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H132").Value = 1552.75
$ws.Range("I132").Value = 1403.6666
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4210.9998
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1680.9998
$ws.Range("N132").Value = -11060
$ws.Range("H135").Value = 51999.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 51999.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 51999.668
$ws.Range("N135").Value = -62139.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1558.1666
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 2187.25
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 6561.75
$ws.Range("M11").Value = -760
$ws.Range("N11").Value = -6841.75
$ws.Range("H12").Value = 84.666664
$ws.Range("I12").Value = 55.75
$ws.Range("J12").Value = 142.5
$ws.Range("K12").Value = 167.25
$ws.Range("L12").Value = 427.5
$ws.Range("M12").Value = 5.75
$ws.Range("N12").Value = -773.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 25000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26640
$ws.Range("H95").Value = 17249.5
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 17249.5
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 17249.5
$ws.Range("N95").Value = -22741.5
$ws.Range("H101").Value = 50000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 50000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H132").Value = 4497
$ws.Range("I132").Value = 4499
$ws.Range("J132").Value = 4495
$ws.Range("K132").Value = 13497
$ws.Range("L132").Value = 13485
$ws.Range("M132").Value = -10967
$ws.Range("N132").Value = -18545
$ws.Range("H135").Value = 39998.25
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 39998.25
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 39998.25
$ws.Range("N135").Value = -50138.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 40000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 40000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 40000
$ws.Range("N74").Value = -41996
$ws.Range("H77").Value = 40000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 40000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 120000
$ws.Range("N77").Value = -129984
$ws.Range("H82").Value = 3275
$ws.Range("I82").Value = 3275
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3275
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2914
$ws.Range("H85").Value = 3275
$ws.Range("I85").Value = 3275
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3275
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2027
$ws.Range("H132").Value = 12955.667
$ws.Range("I132").Value = 13071.333
$ws.Range("J132").Value = 12666.5
$ws.Range("K132").Value = 39213.999
$ws.Range("L132").Value = 37999.5
$ws.Range("M132").Value = -36683.999
$ws.Range("N132").Value = -43059.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 80000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 80000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 80000
$ws.Range("N68").Value = -81622
$ws.Range("H71").Value = 80000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 80000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 240000
$ws.Range("N71").Value = -248112
$ws.Range("H98").Value = 33333
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 33333
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 33333
$ws.Range("N98").Value = -39323
$ws.Range("H101").Value = 20033.666
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 20033.666
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 20033.666
$ws.Range("N101").Value = -26523.666
$ws.Range("H111").Value = 46577.4
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 46577.4
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 46577.4
$ws.Range("N111").Value = -54757.4
$ws.Range("H132").Value = 3073.3333
$ws.Range("I132").Value = 3252.4285
$ws.Range("J132").Value = 2446.5
$ws.Range("K132").Value = 9757.2855
$ws.Range("L132").Value = 7339.5
$ws.Range("M132").Value = -7227.2855
$ws.Range("N132").Value = -12399.5
